$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update row 33 (bug #31): status changes from PENDIENTE to CORREGIDO
# ---------------------------------------------------------------------------
$ws.Range("F33").Value = "CORREGIDO"
$ws.Range("F4").Copy()
$ws.Range("F33").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Row 38 (bug #36) - finish the row that only had the id filled in
# ---------------------------------------------------------------------------
$ws.Range("B38").Value = "Error en el asistente de configuración"
$ws.Range("B22").Copy()
$ws.Range("B38").PasteSpecial(-4122)

$ws.Range("C38").Value = "Configurar edicion: configurar preferencias, elegir 8 equipos, agregar una fase y mostrar el fixture, tocar volver, seleccionar 3 equipos, va a salir el cartel que va a modificarse el fixture, poner aceptar, poner siguiente y se produce un error"
$ws.Range("C24").Copy()
$ws.Range("C38").PasteSpecial(-4122)

$ws.Range("D38").Value = "Tony"

$ws.Range("E38").Value = "admin/edicion/equipos.aspx"
$ws.Range("E38").WrapText = $true

$ws.Range("F38").Value = "CORREGIDO"
$ws.Range("F4").Copy()
$ws.Range("F38").PasteSpecial(-4122)

$ws.Rows.Item(38).RowHeight = 58.5

# ---------------------------------------------------------------------------
# 3. Row 39 (bug #37)
# ---------------------------------------------------------------------------
$ws.Range("A39").Value = 37
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)

$ws.Range("B39").Value = "INTERFAZ PARTIDOS: Arreglos menores"
$ws.Range("B22").Copy()
$ws.Range("B39").PasteSpecial(-4122)

$ws.Range("C39").Value = "Otros partidos de la fecha: colorcito de los estados`n- Widget versus: cambiar foto Partidos Empatados y Partidos Perdidos"
$ws.Range("C24").Copy()
$ws.Range("C39").PasteSpecial(-4122)

$ws.Range("D39").Value = "Facu"

$ws.Range("E39").Value = "torneo/partido"

$ws.Range("F39").Value = "PENDIENTE"
$ws.Range("F18").Copy()
$ws.Range("F39").PasteSpecial(-4122)

$ws.Rows.Item(39).RowHeight = 60

# ---------------------------------------------------------------------------
# 4. Row 40 (bug #38)
# ---------------------------------------------------------------------------
$ws.Range("A40").Value = 38
$ws.Range("A38").Copy()
$ws.Range("A40").PasteSpecial(-4122)

$ws.Range("B40").Value = "Sacar Notificaciones (módulo admin)"
$ws.Range("B22").Copy()
$ws.Range("B40").PasteSpecial(-4122)

$ws.Range("D40").Value = "Facu"

$ws.Range("E40").Value = "admin/"

$ws.Range("F40").Value = "PENDIENTE"
$ws.Range("F18").Copy()
$ws.Range("F40").PasteSpecial(-4122)

$ws.Rows.Item(40).RowHeight = 30

# ---------------------------------------------------------------------------
# 5. Row 41 (bug #39)
# ---------------------------------------------------------------------------
$ws.Range("A41").Value = 39
$ws.Range("A38").Copy()
$ws.Range("A41").PasteSpecial(-4122)

$ws.Range("C41").Value = "Colorcito de los Resultados (Empatado - Perdido - Ganado)`n- Me parece que falta PROXIMOS PARTIDOS!"
$ws.Range("C41").WrapText = $true

$ws.Range("B41").Value = "INTERFAZ EQUIPOS: Arreglos menores"
$ws.Range("B22").Copy()
$ws.Range("B41").PasteSpecial(-4122)

$ws.Range("D41").Value = "Facu"

$ws.Range("E41").Value = "torneo/equipo"

$ws.Range("F41").Value = "PENDIENTE"
$ws.Range("F18").Copy()
$ws.Range("F41").PasteSpecial(-4122)

$ws.Rows.Item(41).RowHeight = 45

# ---------------------------------------------------------------------------
# 6. Rows 42 and 43 - only the id column filled in
# ---------------------------------------------------------------------------
$ws.Range("A42").Value = 40
$ws.Range("A38").Copy()
$ws.Range("A42").PasteSpecial(-4122)

$ws.Range("A43").Value = 41
$ws.Range("A38").Copy()
$ws.Range("A43").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 7. Column E is widened to fit the new, longer "Formulario Asociado" values
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 28

# ---------------------------------------------------------------------------
# 8. Update the view so the selected cell / visible area matches the edit
# ---------------------------------------------------------------------------
$ws.Range("C36").Select()
